$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("`u{FEFF}""page""","short-url","maxPages","items","year","coo_id","coo_name","coo","coo_iso","coa_id","coa_name","coa","coa_iso","refugees","asylum_seekers","returned_refugees","idps","returned_idps","stateless","ooc","oip","hst")

$rowVals = @("1","d8pUGN","1","1","2024","112","Sri Lanka","LKA","LKA","52","Dominica","DMA","DMA","0","5","0","0","0","0","0","-","0")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $rowVals[$i]
}
